$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-30"

# Update the header label cell (shared string "2022 (through 06-29)" -> "2022 (through 06-30)")
$ws.Range("I1").Value = "2022 (through 06-30)"

# Update June's 2022 total (I7): 139 -> 143
$ws.Range("I7").Value = 143

# Update the grand Total row's 2022 total (I14): 802 -> 806
$ws.Range("I14").Value = 806
